$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe.
$ws.Range("D2").Value = "69.308.26"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.854.77"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("D7").Value = "3.850.94"
$ws.Range("E7").Value = "  +2.68%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "4.500.85"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("D16").Value = "3.864.19"
$ws.Range("E16").Value = "  +2.91%  "
$ws.Range("D17").Value = "69.416.22"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("E18").Value = "  +3.83%  "
$ws.Range("E19").Value = "  +7.16%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("E24").Value = "  +4.60%  "
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("D34").Value = "4.004.07"
$ws.Range("E34").Value = "  +2.68%  "
$ws.Range("D35").Value = "3.797.90"
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("E38").Value = "  +4.90%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E48").Value = "  +13.79%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.842.29"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E51").Value = "  +1.18%  "

# Numeric-looking text values (e.g. "603.91", "1.00") must be forced to stay
# text so they are not coerced into numbers by Excel's auto-detection:
# set the cell to Text format, assign, then restore the default "Normal"
# style so no stray number-format style is left behind.
$numericTextCells = @("D5", "D6", "D10", "D11", "D13", "D18", "D19", "D21", "D22", "D25", "D26", "D27", "D28", "D30", "D31", "D33", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D5").Value = "603.91"
$ws.Range("D6").Value = "165.13"
$ws.Range("D10").Value = "0.169"
$ws.Range("D11").Value = "6.34"
$ws.Range("D13").Value = "37.30"
$ws.Range("D18").Value = "7.62"
$ws.Range("D19").Value = "11.74"
$ws.Range("D21").Value = "17.27"
$ws.Range("D22").Value = "490.34"
$ws.Range("D25").Value = "84.71"
$ws.Range("D26").Value = "2.28"
$ws.Range("D27").Value = "12.27"
$ws.Range("D28").Value = "10.10"
$ws.Range("D30").Value = "2.99"
$ws.Range("D31").Value = "8.03"
$ws.Range("D33").Value = "32.49"
$ws.Range("D37").Value = "1.03"
$ws.Range("D38").Value = "0.140"
$ws.Range("D40").Value = "1.00"
$ws.Range("D42").Value = "3.03"
$ws.Range("D43").Value = "440.18"
$ws.Range("D44").Value = "2.00"
$ws.Range("D45").Value = "48.56"
$ws.Range("D46").Value = "8.43"
$ws.Range("D47").Value = "1.00"
$ws.Range("D48").Value = "26.61"
$ws.Range("D49").Value = "143.43"
$ws.Range("D51").Value = "0.0358"

foreach ($ref in $numericTextCells) {
    $ws.Range($ref).Style = "Normal"
}
